$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1733.3334
$ws.Range("I40").Value = 1114.2858
$ws.Range("K40").Value = 1114.2858
$ws.Range("M40").Value = -939.2858000000001
$ws.Range("H74").Value = 8932782
$ws.Range("I74").Value = 3224.5
$ws.Range("K74").Value = 3224.5
$ws.Range("M74").Value = -2288.5
$ws.Range("H77").Value = 8932782
$ws.Range("I77").Value = 3224.5
$ws.Range("K77").Value = 16122.5
$ws.Range("M77").Value = -11442.5
$ws.Range("H113").Value = 4028.5
$ws.Range("I113").Value = 2675
$ws.Range("J113").Value = 5833.1665
$ws.Range("K113").Value = 2675
$ws.Range("L113").Value = 5833.1665
$ws.Range("M113").Value = 579
$ws.Range("N113").Value = -12341.1665
$ws.Range("H129").Value = 1409.5103
$ws.Range("J129").Value = 1448.2128
$ws.Range("L129").Value = 4344.6384
$ws.Range("N129").Value = -14344.6384
$ws.Range("H132").Value = 3086.5557
$ws.Range("I132").Value = 3242.4546
$ws.Range("K132").Value = 9727.363799999999
$ws.Range("M132").Value = -7197.363799999999
$ws.Range("H138").Value = 1509.5682
$ws.Range("J138").Value = 3586.4285
$ws.Range("L138").Value = 10759.2855
$ws.Range("N138").Value = -21039.2855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2064.25
$ws.Range("I2").Value = 1621.5
$ws.Range("J2").Value = 2802.1667
$ws.Range("K2").Value = 1621.5
$ws.Range("L2").Value = 2802.1667
$ws.Range("M2").Value = -1508.5
$ws.Range("N2").Value = -3028.1667
$ws.Range("H76").Value = 14500
$ws.Range("J76").Value = 14500
$ws.Range("L76").Value = 14500
$ws.Range("N76").Value = -15176
$ws.Range("H79").Value = 14500
$ws.Range("J79").Value = 14500
$ws.Range("L79").Value = 14500
$ws.Range("N79").Value = -16840
$ws.Range("H116").Value = 2064.25
$ws.Range("I116").Value = 1621.5
$ws.Range("J116").Value = 2802.1667
$ws.Range("K116").Value = 1621.5
$ws.Range("L116").Value = 2802.1667
$ws.Range("M116").Value = 672.5
$ws.Range("N116").Value = -7390.1667
$ws.Range("H122").Value = 1936.3334
$ws.Range("I122").Value = 1981.5
$ws.Range("J122").Value = 1575
$ws.Range("K122").Value = 5944.5
$ws.Range("L122").Value = 4725
$ws.Range("M122").Value = -3494.5
$ws.Range("N122").Value = -9625
$ws.Range("H132").Value = 19949.572
$ws.Range("I132").Value = 1940.6666
$ws.Range("K132").Value = 5821.9998
$ws.Range("M132").Value = -3291.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2064.25
$ws.Range("I3").Value = 1621.5
$ws.Range("J3").Value = 2802.1667
$ws.Range("K3").Value = 1621.5
$ws.Range("L3").Value = 2802.1667
$ws.Range("M3").Value = -1507.5
$ws.Range("N3").Value = -3030.1667
$ws.Range("H27").Value = 38514
$ws.Range("J27").Value = 38514
$ws.Range("L27").Value = 38514
$ws.Range("N27").Value = -38898
$ws.Range("H105").Value = 4549354.5
$ws.Range("I105").Value = 4950
$ws.Range("J105").Value = 10002640
$ws.Range("K105").Value = 4950
$ws.Range("L105").Value = 10002640
$ws.Range("M105").Value = -3203
$ws.Range("N105").Value = -10006134

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 45000
$ws.Range("J43").Value = 45000
$ws.Range("L43").Value = 45000
$ws.Range("N43").Value = -45368
$ws.Range("H68").Value = 69917.5
$ws.Range("J68").Value = 69917.5
$ws.Range("L68").Value = 69917.5
$ws.Range("N68").Value = -71415.5
$ws.Range("H71").Value = 69917.5
$ws.Range("J71").Value = 69917.5
$ws.Range("L71").Value = 209752.5
$ws.Range("N71").Value = -217240.5
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H88").Value = 34257
$ws.Range("J88").Value = 34257
$ws.Range("L88").Value = 34257
$ws.Range("N88").Value = -35069
$ws.Range("H91").Value = 34257
$ws.Range("J91").Value = 34257
$ws.Range("L91").Value = 34257
$ws.Range("N91").Value = -37065
$ws.Range("H101").Value = 45000
$ws.Range("J101").Value = 45000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -51490

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 20901250
$ws.Range("J37").Value = 20901250
$ws.Range("L37").Value = 62703750
$ws.Range("N37").Value = -62703974
$ws.Range("H131").Value = 734.15
$ws.Range("J131").Value = 738.9286
$ws.Range("L131").Value = 2216.7858
$ws.Range("N131").Value = -12296.7858
$ws.Range("H133").Value = 3238.2222
$ws.Range("J133").Value = 3393
$ws.Range("L133").Value = 10179
$ws.Range("N133").Value = -20299

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4819904
$ws.Range("I70").Value = 24790
$ws.Range("J70").Value = 7816850.5
$ws.Range("K70").Value = 24790
$ws.Range("L70").Value = 7816850.5
$ws.Range("M70").Value = -24520
$ws.Range("N70").Value = -7817390.5
$ws.Range("H73").Value = 4819904
$ws.Range("I73").Value = 24790
$ws.Range("J73").Value = 7816850.5
$ws.Range("K73").Value = 24790
$ws.Range("L73").Value = 7816850.5
$ws.Range("M73").Value = -23854
$ws.Range("N73").Value = -7818722.5
$ws.Range("H80").Value = 3915.3076
$ws.Range("I80").Value = 3566.6667
$ws.Range("J80").Value = 4214.143
$ws.Range("K80").Value = 3566.6667
$ws.Range("L80").Value = 4214.143
$ws.Range("M80").Value = -2568.6667
$ws.Range("N80").Value = -6210.143
$ws.Range("H83").Value = 3915.3076
$ws.Range("I83").Value = 3566.6667
$ws.Range("J83").Value = 4214.143
$ws.Range("K83").Value = 17833.3335
$ws.Range("L83").Value = 21070.715
$ws.Range("M83").Value = -12841.3335
$ws.Range("N83").Value = -31054.715
$ws.Range("H102").Value = 31251098
$ws.Range("I102").Value = 35715396
$ws.Range("K102").Value = 35715396
$ws.Range("M102").Value = -35713774
$ws.Range("H132").Value = 53291.5
$ws.Range("I132").Value = 41855.617
$ws.Range("J132").Value = 127624.75
$ws.Range("K132").Value = 125566.851
$ws.Range("L132").Value = 382874.25
$ws.Range("M132").Value = -123036.851
$ws.Range("N132").Value = -387934.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5900.5
$ws.Range("J22").Value = 6250
$ws.Range("L22").Value = 6250
$ws.Range("N22").Value = -6840
$ws.Range("H27").Value = 5900.5
$ws.Range("J27").Value = 6250
$ws.Range("L27").Value = 6250
$ws.Range("N27").Value = -6464
$ws.Range("H46").Value = 800
$ws.Range("J46").Value = 800
$ws.Range("L46").Value = 800
$ws.Range("N46").Value = -1176
$ws.Range("H129").Value = 28000
$ws.Range("J129").Value = 28000
$ws.Range("L129").Value = 28000
$ws.Range("N129").Value = -38000
$ws.Range("H132").Value = 2717.1428
$ws.Range("I132").Value = 1863.1428
$ws.Range("J132").Value = 3571.1428
$ws.Range("K132").Value = 5589.428400000001
$ws.Range("L132").Value = 10713.4284
$ws.Range("M132").Value = -3059.428400000001
$ws.Range("N132").Value = -15773.4284

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 35389
$ws.Range("J93").Value = 35389
$ws.Range("L93").Value = 35389
$ws.Range("N93").Value = -40381
$ws.Range("H126").Value = 6025
$ws.Range("I126").Value = 20500
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 61500
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -59030
$ws.Range("N126").Value = -8540
